# Auto-generated edit script applying numeric updates to Leviathan_Profits workbook
# Sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 30404.967
$ws.Range("J17").Value = 32457.035
$ws.Range("L17").Value = 97371.105
$ws.Range("N17").Value = -97707.105
$ws.Range("H80").Value = 331.9375
$ws.Range("I80").Value = 212.5
$ws.Range("J80").Value = 403.6
$ws.Range("K80").Value = 637.5
$ws.Range("L80").Value = 1210.8
$ws.Range("M80").Value = 360.5
$ws.Range("N80").Value = -3206.8
$ws.Range("H83").Value = 331.9375
$ws.Range("I83").Value = 212.5
$ws.Range("J83").Value = 403.6
$ws.Range("K83").Value = 1912.5
$ws.Range("L83").Value = 3632.4
$ws.Range("M83").Value = 3079.5
$ws.Range("N83").Value = -13616.4
$ws.Range("H132").Value = 1963.3667
$ws.Range("I132").Value = 1348.2222
$ws.Range("K132").Value = 4044.6666
$ws.Range("M132").Value = -1514.6666
$ws.Range("H135").Value = 1162.3448
$ws.Range("I135").Value = 1010.5
$ws.Range("J135").Value = 1891.2
$ws.Range("K135").Value = 9094.5
$ws.Range("L135").Value = 17020.8
$ws.Range("M135").Value = -6559.5
$ws.Range("N135").Value = -22090.8
$ws.Range("H137").Value = 2672.258
$ws.Range("I137").Value = 2293.2
$ws.Range("K137").Value = 6879.599999999999
$ws.Range("M137").Value = -4329.599999999999
$ws.Range("H138").Value = 2314.806
$ws.Range("J138").Value = 3458.6667
$ws.Range("L138").Value = 10376.0001
$ws.Range("N138").Value = -20656.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 4874.7144
$ws.Range("J22").Value = 6229.5
$ws.Range("L22").Value = 6229.5
$ws.Range("N22").Value = -6827.5
$ws.Range("H32").Value = 43287.35
$ws.Range("I32").Value = 26438.324
$ws.Range("K32").Value = 26438.324
$ws.Range("M32").Value = -26151.324
$ws.Range("H45").Value = 483990.75
$ws.Range("I45").Value = 596774.75
$ws.Range("K45").Value = 596774.75
$ws.Range("M45").Value = -596397.75
$ws.Range("H132").Value = 11626.9
$ws.Range("I132").Value = 12180.947
$ws.Range("K132").Value = 36542.841
$ws.Range("M132").Value = -34012.841

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H75").Value = 6499
$ws.Range("I75").Value = 6499
$ws.Range("K75").Value = 6499
$ws.Range("M75").Value = -5563
$ws.Range("H78").Value = 6499
$ws.Range("I78").Value = 6499
$ws.Range("K78").Value = 19497
$ws.Range("M78").Value = -14817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2702.5386
$ws.Range("I31").Value = 2093.25
$ws.Range("K31").Value = 2093.25
$ws.Range("M31").Value = -1798.25
$ws.Range("H34").Value = 2702.5386
$ws.Range("I34").Value = 2093.25
$ws.Range("K34").Value = 2093.25
$ws.Range("M34").Value = -1891.25
$ws.Range("H35").Value = 2666.6667
$ws.Range("I35").Value = 2500
$ws.Range("K35").Value = 2500
$ws.Range("M35").Value = -2206
$ws.Range("H58").Value = 1300.0938
$ws.Range("I58").Value = 1206.8572
$ws.Range("J58").Value = 1478.091
$ws.Range("K58").Value = 1206.8572
$ws.Range("L58").Value = 1478.091
$ws.Range("M58").Value = -1003.8572
$ws.Range("N58").Value = -1884.091
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H103").Value = 19887
$ws.Range("I103").Value = 19887
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 19887
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -18715
$ws.Range("N103").ClearContents()
$ws.Range("H122").Value = 4605.625
$ws.Range("I122").Value = 3440.8333
$ws.Range("K122").Value = 10322.4999
$ws.Range("M122").Value = -7872.499899999999
$ws.Range("H132").Value = 2380.7334
$ws.Range("I132").Value = 2324
$ws.Range("J132").Value = 2749.5
$ws.Range("K132").Value = 6972
$ws.Range("L132").Value = 8248.5
$ws.Range("M132").Value = -4442
$ws.Range("N132").Value = -13308.5
$ws.Range("H134").Value = 2480.606
$ws.Range("I134").Value = 2595.36
$ws.Range("J134").Value = 2122
$ws.Range("K134").Value = 7786.08
$ws.Range("L134").Value = 6366
$ws.Range("M134").Value = -5251.08
$ws.Range("N134").Value = -11436
$ws.Range("H136").Value = 1300.0938
$ws.Range("I136").Value = 1206.8572
$ws.Range("J136").Value = 1478.091
$ws.Range("K136").Value = 3620.5716
$ws.Range("L136").Value = 4434.272999999999
$ws.Range("M136").Value = -1070.5716
$ws.Range("N136").Value = -9534.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 207145.2
$ws.Range("I48").Value = 503874
$ws.Range("J48").Value = 9326
$ws.Range("K48").Value = 1511622
$ws.Range("L48").Value = 27978
$ws.Range("M48").Value = -1511372
$ws.Range("N48").Value = -28478

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -6058
$ws.Range("H126").Value = 3908
$ws.Range("I126").Value = 3825.375
$ws.Range("J126").Value = 4128.3335
$ws.Range("K126").Value = 11476.125
$ws.Range("L126").Value = 12385.0005
$ws.Range("M126").Value = -9006.125
$ws.Range("N126").Value = -17325.0005
$ws.Range("H132").Value = 3017.8462
$ws.Range("I132").Value = 3085
$ws.Range("J132").Value = 2835.5715
$ws.Range("K132").Value = 9255
$ws.Range("L132").Value = 8506.7145
$ws.Range("M132").Value = -6725
$ws.Range("N132").Value = -13566.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18299.277
$ws.Range("J7").Value = 6111.2856
$ws.Range("L7").Value = 6111.2856
$ws.Range("N7").Value = -6335.2856
$ws.Range("H31").Value = 2110.2856
$ws.Range("I31").Value = 2135.818
$ws.Range("J31").Value = 2016.6666
$ws.Range("K31").Value = 2135.818
$ws.Range("L31").Value = 2016.6666
$ws.Range("M31").Value = -1887.818
$ws.Range("N31").Value = -2512.6666
$ws.Range("H40").Value = 6717.923
$ws.Range("I40").Value = 6734.8
$ws.Range("K40").Value = 6734.8
$ws.Range("M40").Value = -6598.8
$ws.Range("H46").Value = 108905.75
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376
$ws.Range("H61").Value = 68334.1
$ws.Range("I61").Value = 68344.07
$ws.Range("J61").Value = 68304.2
$ws.Range("K61").Value = 68344.07
$ws.Range("L61").Value = 68304.2
$ws.Range("M61").Value = -68142.07
$ws.Range("N61").Value = -68708.2
$ws.Range("H80").Value = 24752
$ws.Range("J80").Value = 24752
$ws.Range("L80").Value = 24752
$ws.Range("N80").Value = -26998
$ws.Range("H83").Value = 24752
$ws.Range("J83").Value = 24752
$ws.Range("L83").Value = 74256
$ws.Range("N83").Value = -85488
$ws.Range("H113").Value = 68334.1
$ws.Range("I113").Value = 68344.07
$ws.Range("J113").Value = 68304.2
$ws.Range("K113").Value = 68344.07
$ws.Range("L113").Value = 68304.2
$ws.Range("M113").Value = -66174.07
$ws.Range("N113").Value = -72644.2
$ws.Range("H125").Value = 69715
$ws.Range("J125").Value = 69715
$ws.Range("L125").Value = 69715
$ws.Range("N125").Value = -79555
$ws.Range("H126").Value = 18299.277
$ws.Range("J126").Value = 6111.2856
$ws.Range("L126").Value = 18333.8568
$ws.Range("N126").Value = -23273.8568
$ws.Range("H132").Value = 4674.222
$ws.Range("I132").Value = 4441
$ws.Range("J132").Value = 5140.6665
$ws.Range("K132").Value = 13323
$ws.Range("L132").Value = 15421.9995
$ws.Range("M132").Value = -10793
$ws.Range("N132").Value = -20481.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4318.2593
$ws.Range("J81").Value = 4000.75
$ws.Range("L81").Value = 8001.5
$ws.Range("N81").Value = -10123.5
$ws.Range("H84").Value = 4318.2593
$ws.Range("J84").Value = 4000.75
$ws.Range("L84").Value = 40007.5
$ws.Range("N84").Value = -50615.5
$ws.Range("H100").Value = 3900.125
$ws.Range("I100").Value = 2693.4666
$ws.Range("K100").Value = 5386.9332
$ws.Range("M100").Value = -4845.9332
$ws.Range("H136").Value = 1017.4286
$ws.Range("I136").Value = 824.4
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2473.2
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 76.80000000000018
$ws.Range("N136").Value = -9600

